# full load logic - weights and scores
#
# Fill in the "Notes" column (C) for the Packaging section of the
# "Product" sheet: the "Easy to open" requirement gets a note about
# getting cut, and the "Recyclable" requirement gets a note about toxic
# sludge. New cells copy the formatting already used by their row's
# "Score" cell (column B) so they pick up the same border/font as the
# rest of the scored rows instead of the blank placeholder formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "Cut my finger"

$ws.Range("B14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Made from toxic sludge"

# Widen the new Notes column so the text isn't clipped.
$ws.Columns.Item(3).ColumnWidth = 19.29
